# Refresh the crypto price/volume table with the latest scrape.
# (GitHub Actions job: "Updated cryptos list")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "25.942.03"  # D2 Price
$ws.Cells.Item(2, 5).Value = "  +0.44%  "  # E2 Volume(1h)
$ws.Cells.Item(3, 4).Value = "1.647.02"  # D3 Price
$ws.Cells.Item(3, 5).Value = "  +0.85%  "  # E3 Volume(1h)
$ws.Cells.Item(4, 4).Value = "'1.007"  # D4 Price
$ws.Cells.Item(4, 5).Value = "  +0.58%  "  # E4 Volume(1h)
$ws.Cells.Item(5, 4).Value = "'216.11"  # D5 Price
$ws.Cells.Item(5, 5).Value = "  +0.67%  "  # E5 Volume(1h)
$ws.Cells.Item(6, 4).Value = "'0.5102"  # D6 Price
$ws.Cells.Item(6, 5).Value = "  +1.58%  "  # E6 Volume(1h)
$ws.Cells.Item(7, 4).Value = "'1.007"  # D7 Price
$ws.Cells.Item(7, 5).Value = "  +0.54%  "  # E7 Volume(1h)
$ws.Cells.Item(8, 4).Value = "'0.2579"  # D8 Price
$ws.Cells.Item(8, 5).Value = "  +0.57%  "  # E8 Volume(1h)
$ws.Cells.Item(9, 4).Value = "'0.06424"  # D9 Price
$ws.Cells.Item(9, 5).Value = "  +0.66%  "  # E9 Volume(1h)
$ws.Cells.Item(10, 4).Value = "'19.67"  # D10 Price
$ws.Cells.Item(10, 5).Value = "  +0.13%  "  # E10 Volume(1h)
$ws.Cells.Item(11, 4).Value = "'0.07783"  # D11 Price
$ws.Cells.Item(11, 5).Value = "  +1.15%  "  # E11 Volume(1h)
$ws.Cells.Item(12, 4).Value = "'4.324"  # D12 Price
$ws.Cells.Item(12, 5).Value = "  +1.94%  "  # E12 Volume(1h)
$ws.Cells.Item(13, 4).Value = "1.651.41"  # D13 Price
$ws.Cells.Item(13, 5).Value = "  +1.09%  "  # E13 Volume(1h)
$ws.Cells.Item(14, 4).Value = "'0.5467"  # D14 Price
$ws.Cells.Item(14, 5).Value = "  +0.85%  "  # E14 Volume(1h)
$ws.Cells.Item(15, 4).Value = "0.0₅7895"  # D15 Price
$ws.Cells.Item(15, 5).Value = "  -0.24%  "  # E15 Volume(1h)
$ws.Cells.Item(16, 4).Value = "'64.75"  # D16 Price
$ws.Cells.Item(16, 5).Value = "  +1.96%  "  # E16 Volume(1h)
$ws.Cells.Item(17, 4).Value = "26.032.28"  # D17 Price
$ws.Cells.Item(17, 5).Value = "  +0.78%  "  # E17 Volume(1h)
$ws.Cells.Item(18, 4).Value = "'1.007"  # D18 Price
$ws.Cells.Item(18, 5).Value = "  +0.54%  "  # E18 Volume(1h)
$ws.Cells.Item(19, 4).Value = "'198.51"  # D19 Price
$ws.Cells.Item(19, 5).Value = "  -1.66%  "  # E19 Volume(1h)
$ws.Cells.Item(20, 4).Value = "'4.481"  # D20 Price
$ws.Cells.Item(20, 5).Value = "  +3.57%  "  # E20 Volume(1h)
$ws.Cells.Item(21, 4).Value = "'10.02"  # D21 Price
$ws.Cells.Item(21, 5).Value = "  +1.03%  "  # E21 Volume(1h)
$ws.Cells.Item(22, 4).Value = "'6.074"  # D22 Price
$ws.Cells.Item(22, 5).Value = "  +1.96%  "  # E22 Volume(1h)
$ws.Cells.Item(23, 4).Value = "'1.010"  # D23 Price
$ws.Cells.Item(23, 5).Value = "  +0.74%  "  # E23 Volume(1h)
$ws.Cells.Item(24, 4).Value = "'1.863"  # D24 Price
$ws.Cells.Item(24, 5).Value = "  -3.69%  "  # E24 Volume(1h)
$ws.Cells.Item(25, 4).Value = "'140.35"  # D25 Price
$ws.Cells.Item(25, 5).Value = "  -0.89%  "  # E25 Volume(1h)
$ws.Cells.Item(26, 4).Value = "'0.1151"  # D26 Price
$ws.Cells.Item(26, 5).Value = "  +0.93%  "  # E26 Volume(1h)
$ws.Cells.Item(27, 4).Value = "'6.906"  # D27 Price
$ws.Cells.Item(27, 5).Value = "  +3.05%  "  # E27 Volume(1h)
$ws.Cells.Item(28, 4).Value = "'15.76"  # D28 Price
$ws.Cells.Item(28, 5).Value = "  +0.59%  "  # E28 Volume(1h)
$ws.Cells.Item(29, 4).Value = "'1.243"  # D29 Price
$ws.Cells.Item(29, 5).Value = "  +0.44%  "  # E29 Volume(1h)
$ws.Cells.Item(30, 4).Value = "'0.05022"  # D30 Price
$ws.Cells.Item(30, 5).Value = "  +0.44%  "  # E30 Volume(1h)
$ws.Cells.Item(31, 4).Value = "'3.288"  # D31 Price
$ws.Cells.Item(31, 5).Value = "  +0.83%  "  # E31 Volume(1h)
$ws.Cells.Item(32, 4).Value = "'3.203"  # D32 Price
$ws.Cells.Item(32, 5).Value = "  +0.78%  "  # E32 Volume(1h)
$ws.Cells.Item(33, 4).Value = "'1.547"  # D33 Price
$ws.Cells.Item(33, 5).Value = "  +0.54%  "  # E33 Volume(1h)
$ws.Cells.Item(34, 4).Value = "'2.366"  # D34 Price
$ws.Cells.Item(34, 5).Value = "  +0.03%  "  # E34 Volume(1h)
$ws.Cells.Item(35, 4).Value = "'0.8950"  # D35 Price
$ws.Cells.Item(35, 5).Value = "  +0.47%  "  # E35 Volume(1h)
$ws.Cells.Item(36, 4).Value = "'2.596"  # D36 Price
$ws.Cells.Item(36, 5).Value = "  -0.69%  "  # E36 Volume(1h)
$ws.Cells.Item(37, 4).Value = "1.135.62"  # D37 Price
$ws.Cells.Item(37, 5).Value = "  -3.05%  "  # E37 Volume(1h)
$ws.Cells.Item(38, 4).Value = "'0.5543"  # D38 Price
$ws.Cells.Item(38, 5).Value = "  -0.77%  "  # E38 Volume(1h)
$ws.Cells.Item(39, 4).Value = "'0.01567"  # D39 Price
$ws.Cells.Item(39, 5).Value = "  +0.41%  "  # E39 Volume(1h)
$ws.Cells.Item(40, 4).Value = "'1.008"  # D40 Price
$ws.Cells.Item(40, 5).Value = "  +0.73%  "  # E40 Volume(1h)
$ws.Cells.Item(44, 4).Value = "'100.05"  # D44 Price
$ws.Cells.Item(44, 5).Value = "  +0.72%  "  # E44 Volume(1h)

# Rows 41-51: one new coin (mCoin) entered the top-50 list, pushing the
# remaining rows down by one and dropping Algorand off the bottom.
$ws.Cells.Item(41, 2).Value = "mCoin"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Cells.Item(41, 4).Value = "'2.551"
$ws.Cells.Item(41, 5).Value = "  -0.58%  "
$ws.Cells.Item(42, 2).Value = "FraxShare"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(42, 4).Value = "'5.666"
$ws.Cells.Item(42, 5).Value = "  -0.31%  "
$ws.Cells.Item(43, 2).Value = "TrustWalletToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(43, 4).Value = "'0.8174"
$ws.Cells.Item(43, 5).Value = "  +1.53%  "
$ws.Cells.Item(45, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(45, 4).Value = "0.0₈125"
$ws.Cells.Item(45, 5).Value = "  +9.17%  "
$ws.Cells.Item(46, 2).Value = "RocketPoolETH"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(46, 4).Value = "1.787.22"
$ws.Cells.Item(46, 5).Value = "  +0.99%  "
$ws.Cells.Item(47, 2).Value = "Mantle"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(47, 4).Value = "'0.4537"
$ws.Cells.Item(47, 5).Value = "  +0.55%  "
$ws.Cells.Item(49, 2).Value = "Aave"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(49, 4).Value = "'55.25"
$ws.Cells.Item(49, 5).Value = "  +1.00%  "
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).Value = "'0.05091"
$ws.Cells.Item(50, 5).Value = "  +0.25%  "
$ws.Cells.Item(51, 2).Value = "USDD"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Cells.Item(51, 4).Value = "'1.008"
$ws.Cells.Item(51, 5).Value = "  +0.70%  "
